$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I16").Value = "ОК"
